# fix: fixed formatting when scrapping floating point numbers
#
# "Importe" (column H) values were scraped in es-AR style, e.g. "1.234,56"
# (dot thousands separator, comma decimal separator). Re-store each as a
# plain decimal-point string, e.g. "1234.56", matching the fixed scraper
# output. A couple of "Razon social" names (column E) also had a stray
# comma where a period belonged; fix those too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$importeUpdates = @(
  ,@(2, "405.00")
  ,@(3, "95.00")
  ,@(4, "1600.00")
  ,@(5, "1870.00")
  ,@(6, "29749.99")
  ,@(7, "86346.21")
  ,@(8, "15.00")
  ,@(9, "597.00")
  ,@(10, "1271.00")
  ,@(11, "821.00")
  ,@(12, "16170.89")
  ,@(13, "51.12")
  ,@(14, "9376.96")
  ,@(15, "3100.00")
  ,@(16, "260.00")
  ,@(17, "40.00")
  ,@(18, "28184.00")
  ,@(19, "146.00")
  ,@(20, "511.80")
  ,@(21, "36186.00")
  ,@(22, "4908.85")
  ,@(23, "42830.00")
  ,@(24, "591.25")
  ,@(25, "9074.50")
  ,@(26, "90.80")
  ,@(27, "7337.20")
  ,@(28, "33920.70")
  ,@(29, "1097.00")
  ,@(30, "2377.00")
  ,@(31, "11752.69")
  ,@(32, "134.00")
  ,@(33, "13390.98")
  ,@(34, "142.00")
  ,@(35, "52.20")
  ,@(36, "14.00")
  ,@(37, "218.14")
  ,@(38, "336.00")
  ,@(39, "4485.00")
  ,@(40, "22760.00")
  ,@(41, "18295.00")
  ,@(42, "308.69")
  ,@(43, "6750.00")
  ,@(44, "4200.00")
  ,@(45, "1972.30")
  ,@(46, "35011.75")
  ,@(47, "4008.60")
  ,@(48, "665.50")
  ,@(49, "3452.60")
  ,@(50, "118.40")
  ,@(51, "52152.00")
  ,@(52, "26470.00")
  ,@(53, "7715.90")
  ,@(54, "35145.50")
  ,@(55, "850.00")
  ,@(56, "9212.00")
  ,@(57, "525.00")
  ,@(58, "957.00")
  ,@(59, "841.00")
  ,@(60, "202.87")
  ,@(61, "286.00")
  ,@(62, "3293.62")
  ,@(63, "200.00")
  ,@(64, "107070.00")
  ,@(65, "1.16")
  ,@(66, "110.96")
  ,@(67, "7788.85")
  ,@(68, "1874.06")
  ,@(69, "2652.32")
  ,@(70, "8218.84")
  ,@(71, "1741.60")
  ,@(72, "6941.25")
  ,@(73, "7166.00")
  ,@(74, "265.05")
  ,@(75, "165.50")
  ,@(76, "198.00")
  ,@(77, "72.00")
  ,@(78, "3790.00")
  ,@(79, "412.00")
  ,@(80, "7195.00")
  ,@(81, "1012.12")
  ,@(82, "17500.00")
  ,@(83, "953.48")
  ,@(84, "2250.00")
  ,@(85, "24000.00")
  ,@(86, "175.00")
  ,@(87, "5920.00")
  ,@(88, "2465.00")
  ,@(89, "2130.00")
  ,@(90, "125000.00")
  ,@(91, "1200.00")
  ,@(92, "140.00")
  ,@(93, "40.00")
  ,@(94, "1185000.00")
  ,@(95, "30557.43")
  ,@(96, "250.00")
  ,@(97, "1800.00")
  ,@(98, "45.00")
  ,@(99, "9265.70")
  ,@(100, "26983.00")
  ,@(101, "250.00")
  ,@(102, "16000.00")
  ,@(103, "1210.00")
  ,@(104, "9171.52")
  ,@(105, "750.00")
  ,@(106, "900.00")
  ,@(107, "7719.80")
  ,@(108, "3240.00")
  ,@(109, "3000.00")
  ,@(110, "400.00")
  ,@(111, "2500.00")
  ,@(112, "16780.00")
  ,@(113, "625.00")
  ,@(114, "175.00")
  ,@(115, "153.36")
  ,@(116, "245.00")
  ,@(117, "185.69")
  ,@(118, "175.45")
  ,@(119, "4456.00")
  ,@(120, "3449.00")
  ,@(121, "75.00")
  ,@(122, "224.00")
  ,@(123, "258.34")
  ,@(124, "638.10")
  ,@(125, "2238.28")
  ,@(126, "158.40")
  ,@(127, "21157.20")
  ,@(128, "33032.40")
  ,@(129, "3240.00")
  ,@(130, "3234.83")
  ,@(131, "210.00")
  ,@(132, "49950.00")
  ,@(133, "823153.86")
  ,@(134, "278500.00")
  ,@(135, "17000.00")
  ,@(136, "157000.00")
  ,@(137, "18480.00")
  ,@(138, "151829.00")
  ,@(139, "212500.00")
  ,@(140, "294000.00")
  ,@(141, "266000.00")
  ,@(142, "147500.00")
  ,@(143, "6700.00")
  ,@(144, "7902.00")
)

foreach ($pair in $importeUpdates) {
  $r = $pair[0]
  $newVal = $pair[1]
  $cell = $ws.Cells.Item($r, 8)
  # Leading apostrophe forces text (quote-prefixed), keeping the cell a
  # shared-string / text cell instead of letting Excel coerce it into a
  # real number (which would silently drop the trailing zeros, e.g.
  # "405.00" -> 405).
  $cell.Value = "'" + $newVal
}

$nameFixes = @{
  "IZAGUIRRE CARLOS MARIA, MOREND MARIA ELENA Y MOREND MARIA TERESA" = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
  "FERNANDEZ MARIO H, GALLICET OSCAR M" = "FERNANDEZ MARIO H. GALLICET OSCAR M"
}

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
  $cell = $ws.Cells.Item($r, 5)
  $v = $cell.Value2
  if ($nameFixes.ContainsKey($v)) {
    $cell.Value = $nameFixes[$v]
  }
}
